# ActivePlayers_BAN.xlsx - additional scraping update
# 1) Bump a handful of existing TEST/ODI/T20 counters (new match data scraped).
# 2) Insert a brand-new player row (Rishad Hossain) before the current "Rony
#    Talukdar" row, pushing Rony Talukdar .. Zakir Hasan down by one row.
# 3) Bump a few counters on rows that got shifted down as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: value bumps above the insertion point (rows 1-26 untouched in shape) ---
$ws.Cells.Item(6, 3).Value  = 19   # Ebadat Hossain      TEST 18 -> 19
$ws.Cells.Item(7, 5).Value  = 16   # Hasan Mahmud        T20  14 -> 16
$ws.Cells.Item(8, 3).Value  = 12   # Khaled Ahmed        TEST 11 -> 12
$ws.Cells.Item(9, 3).Value  = 38   # Litton Das          TEST 37 -> 38
$ws.Cells.Item(9, 5).Value  = 71   # Litton Das          T20  69 -> 71
$ws.Cells.Item(13, 3).Value = 38   # Mehidy Hasan Miraz  TEST 37 -> 38
$ws.Cells.Item(13, 5).Value = 23   # Mehidy Hasan Miraz  T20  22 -> 23
$ws.Cells.Item(17, 3).Value = 56   # Mominul Haque       TEST 55 -> 56
$ws.Cells.Item(20, 3).Value = 85   # Mushfiqur Rahim     TEST 84 -> 85
$ws.Cells.Item(21, 5).Value = 83   # Mustafizur Rahman   T20  82 -> 83
$ws.Cells.Item(22, 3).Value = 22   # Najmul Hossain Shanto TEST 21 -> 22
$ws.Cells.Item(22, 5).Value = 23   # Najmul Hossain Shanto T20  21 -> 23
$ws.Cells.Item(23, 5).Value = 33   # Nasum Ahmed         T20  31 -> 33

# --- Step 2: insert the new player row at row 27, shifting everything after it down ---
$ws.Rows(27).Insert()

$ws.Cells.Item(27, 1).Value = "Rishad Hossain"
$ws.Cells.Item(27, 2).Value = "'7198"
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = 1

# --- Step 3: value bumps on rows that shifted down (now rows 28-44) ---
$ws.Cells.Item(28, 5).Value = 7    # Rony Talukdar       T20  5 -> 7
$ws.Cells.Item(33, 3).Value = 66   # Shakib Al Hasan     TEST 65 -> 66
$ws.Cells.Item(33, 5).Value = 115  # Shakib Al Hasan     T20  113 -> 115
$ws.Cells.Item(34, 5).Value = 15   # Shamim Hossain      T20  13 -> 15
$ws.Cells.Item(36, 3).Value = 6    # Shoriful Islam      TEST 5 -> 6
$ws.Cells.Item(36, 5).Value = 30   # Shoriful Islam      T20  29 -> 30
$ws.Cells.Item(38, 3).Value = 41   # Taijul Islam        TEST 40 -> 41
$ws.Cells.Item(39, 3).Value = 70   # Tamim Iqbal         TEST 69 -> 70
$ws.Cells.Item(41, 5).Value = 52   # Taskin Ahmed        T20  50 -> 52
$ws.Cells.Item(42, 5).Value = 6    # Towhid Hridoy       T20  4 -> 6
